$wb = $excel.ActiveWorkbook

# Sheet "展览" - update "想去人数" (want-to-go count) column F
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 275
$ws1.Range("F3").Value = 94
$ws1.Range("F4").Value = 1043
$ws1.Range("F5").Value = 558

# Sheet "全部类型" - update "想去人数" (want-to-go count) column F
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 275
$ws4.Range("F3").Value = 94
$ws4.Range("F4").Value = 1043
$ws4.Range("F6").Value = 558
